# Update E2/F2 and E3/F3 values on the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("E2").Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Range("F2").Value = 680

    $ws.Range("E3").Value = "2024.05.01 09:00-05.01 17:00"
    $ws.Range("F3").Value = 57
}
